# Auto-generated Excel COM-interop script
# Applies the cell value updates described by the provided diff
# (Jogos_da_Semana_FlashScore_2024-10-18.xlsx odds update)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("G9").Value = 2.1
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.6
$ws.Range("J9").Value = 2.88
$ws.Range("K9").Value = 2
$ws.Range("N9").Value = 7.5
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.35
$ws.Range("R9").Value = 1.57
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("AA9").Value = 19
$ws.Range("AC9").Value = 7.5
$ws.Range("AD9").Value = 6.5
$ws.Range("AG9").Value = 451
$ws.Range("AH9").Value = 8.5
$ws.Range("AO9").Value = 12
$ws.Range("AS9").Value = 201
$ws.Range("AT9").Value = 2.5
$ws.Range("AU9").Value = 8.5

# Row 25
$ws.Range("J25").Value = 2.5
$ws.Range("K25").Value = 2.2
$ws.Range("M25").Value = 1.05
$ws.Range("N25").Value = 11
$ws.Range("O25").Value = 1.3
$ws.Range("P25").Value = 3.5
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 1.9
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("X25").Value = 9
$ws.Range("AB25").Value = 26
$ws.Range("AC25").Value = 10
$ws.Range("AG25").Value = 251
$ws.Range("AN25").Value = 4
$ws.Range("AX25").Value = 21

# Row 27
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 2.88
$ws.Range("I27").Value = 4.5
$ws.Range("J27").Value = 2.88
$ws.Range("K27").Value = 1.91
$ws.Range("M27").Value = 1.13
$ws.Range("N27").Value = 6
$ws.Range("O27").Value = 1.5
$ws.Range("P27").Value = 2.5
$ws.Range("Q27").Value = 2.6
$ws.Range("S27").Value = 1.57
$ws.Range("T27").Value = 2.25
$ws.Range("X27").Value = 8
$ws.Range("Z27").Value = 17
$ws.Range("AC27").Value = 6
$ws.Range("AH27").Value = 9
$ws.Range("AI27").Value = 21
$ws.Range("AJ27").Value = 17
$ws.Range("AK27").Value = 51
$ws.Range("AS27").Value = 251
$ws.Range("AT27").Value = 2.25

# Row 34
$ws.Range("K34").Value = 1.95
$ws.Range("O34").Value = 1.44
$ws.Range("P34").Value = 2.63
$ws.Range("R34").Value = 1.5

# Row 62
$ws.Range("G62").Value = 2.05
$ws.Range("H62").Value = 3.1
$ws.Range("I62").Value = 4.2
$ws.Range("J62").Value = 2.75
$ws.Range("L62").Value = 4.5
$ws.Range("M62").Value = 1.05
$ws.Range("O62").Value = 1.37
$ws.Range("Q62").Value = 2.25
$ws.Range("R62").Value = 1.62
$ws.Range("U62").Value = 1.95
$ws.Range("V62").Value = 1.8
$ws.Range("W62").Value = 6.5
$ws.Range("X62").Value = 9
$ws.Range("Y62").Value = 9
$ws.Range("Z62").Value = 17
$ws.Range("AD62").Value = 6
$ws.Range("AG62").Value = 351
$ws.Range("AH62").Value = 10
$ws.Range("AI62").Value = 19
$ws.Range("AJ62").Value = 15
$ws.Range("AL62").Value = 34
$ws.Range("AO62").Value = 12
$ws.Range("AP62").Value = 23
$ws.Range("AW62").Value = 5.5
$ws.Range("AX62").Value = 21
$ws.Range("AY62").Value = 34
$ws.Range("AZ62").Value = 81
$ws.Range("BB62").Value = 301

# Row 63
$ws.Range("H63").Value = 3.9

# Row 64
$ws.Range("G64").Value = 1.5
$ws.Range("H64").Value = 4.15
$ws.Range("I64").Value = 5.6
$ws.Range("J64").Value = 1.95
$ws.Range("K64").Value = 2.4
$ws.Range("L64").Value = 5.2
$ws.Range("M64").Value = 1.02
$ws.Range("N64").Value = 14.4
$ws.Range("Q64").Value = 1.5
$ws.Range("R64").Value = 2.27
$ws.Range("S64").Value = 1.26
$ws.Range("T64").Value = 3.62
$ws.Range("U64").Value = 1.6
$ws.Range("V64").Value = 2.07
$ws.Range("W64").Value = 9
$ws.Range("X64").Value = 8.5
$ws.Range("AA64").Value = 11
$ws.Range("AB64").Value = 19.5
$ws.Range("AC64").Value = 15.5
$ws.Range("AD64").Value = 8.5
$ws.Range("AE64").Value = 14.5
$ws.Range("AF64").Value = 50
$ws.Range("AG64").Value = 300
$ws.Range("AH64").Value = 19.5
$ws.Range("AI64").Value = 40
$ws.Range("AJ64").Value = 17.5
$ws.Range("AK64").Value = 110
$ws.Range("AL64").Value = 50
$ws.Range("AM64").Value = 45
$ws.Range("AN64").Value = 3.55
$ws.Range("AO64").Value = 6.9
$ws.Range("AP64").Value = 13.5
$ws.Range("AQ64").Value = 19
$ws.Range("AR64").Value = 37
$ws.Range("AT64").Value = 3.35
$ws.Range("AU64").Value = 7
$ws.Range("AV64").Value = 50
$ws.Range("AW64").Value = 7.3
$ws.Range("AX64").Value = 29
$ws.Range("AZ64").Value = 175
$ws.Range("BA64").Value = 175

# Row 81
$ws.Range("M81").Value = 1.02
$ws.Range("N81").Value = 17
$ws.Range("O81").Value = 1.13
$ws.Range("P81").Value = 5
$ws.Range("Q81").Value = 1.57
$ws.Range("R81").Value = 2.35
$ws.Range("S81").Value = 1.3
